$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old date number format / column width before writing new values
$ws.Range("A1:F1").ClearFormats()
$ws.Columns.Item(2).ClearFormats()

$ws.Range("A1").Value = "Dispersion in X column: "
$ws.Range("B1").Value = 0.30020066935891704
$ws.Range("C1").Value = "Dispersion in Y column: "
$ws.Range("D1").Value = 0.053782220004107456
$ws.Range("E1").Value = "Dispersion in Z column: "
$ws.Range("F1").Value = 7.427918040272801
